$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin / Link / Volume(1h) text cells
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("E26").Value = "  +4.08%  "
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("E29").Value = "  +4.91%  "
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("E38").Value = "  +4.73%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("E40").Value = "  +3.17%  "
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E51").Value = "  +2.62%  "

# Update Price column cells, forcing text so values like "21.98" are not
# auto-converted to numbers by Excel (column stores text-formatted prices).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"
$ws.Range("D2").Value = "26.293.67"
$ws.Range("D3").Value = "1.679.34"
$ws.Range("D5").Value = "218.18"
$ws.Range("D6").Value = "0.5250"
$ws.Range("D8").Value = "0.2693"
$ws.Range("D9").Value = "0.06468"
$ws.Range("D10").Value = "21.98"
$ws.Range("D11").Value = "0.07528"
$ws.Range("D12").Value = "4.530"
$ws.Range("D13").Value = "1.658.38"
$ws.Range("D14").Value = "0.5802"
$ws.Range("D15").Value = "0.000008516"
$ws.Range("D16").Value = "64.74"
$ws.Range("D17").Value = "26.334.25"
$ws.Range("D18").Value = "4.921"
$ws.Range("D20").Value = "10.87"
$ws.Range("D21").Value = "190.02"
$ws.Range("D22").Value = "6.208"
$ws.Range("D23").Value = "1.007"
$ws.Range("D24").Value = "145.56"
$ws.Range("D25").Value = "7.818"
$ws.Range("D27").Value = "15.80"
$ws.Range("D28").Value = "0.06460"
$ws.Range("D30").Value = "1.326"
$ws.Range("D31").Value = "3.594"
$ws.Range("D32").Value = "3.603"
$ws.Range("D34").Value = "1.029"
$ws.Range("D35").Value = "0.6241"
$ws.Range("D36").Value = "2.404"
$ws.Range("D37").Value = "2.733"
$ws.Range("D38").Value = "6.440"
$ws.Range("D39").Value = "0.01625"
$ws.Range("D40").Value = "1.108.56"
$ws.Range("D41").Value = "0.8767"
$ws.Range("D43").Value = "100.68"
$ws.Range("D44").Value = "1.831.91"
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("D46").Value = "56.94"
$ws.Range("D47").Value = "8.209"
$ws.Range("D48").Value = "1.006"
$ws.Range("D49").Value = "0.05268"
$ws.Range("D50").Value = "0.4289"
$ws.Range("D51").Value = "6.083"
$priceRange.Style = "Normal"

Write-Host "Updated cryptos list"
